# Implemented v1 flow for setup of new player
# - Remove the "Explorer" row from the Cards table (old row 47)
# - Update the quantity for "Scout" (now row 47) from 8 to 5
# - Apply an explicit General number format to the "Id" cell of the new
#   last row (L48, formerly the "Viper" row's Id cell)
# - Update the current sheet selection / scroll position

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete the "Explorer" row (row 47). Excel shifts rows 48-49 up to 47-48
# and shrinks the Table1 structured range + sheet dimension accordingly.
$ws.Rows.Item(47).Delete()

# The old "Scout" row (previously row 48) is now row 47; its Quantity
# changes from 8 to 5.
$ws.Range("A47").Value = 5

# The old "Viper" row (previously row 49) is now row 48; its Id cell
# (L48) gets an explicit "General" number format applied.
$ws.Range("L48").NumberFormat = "General"

# Update the view: scroll so row 22 is at the top and select B49 (the
# cell below the table, where a new row would be typed next).
$win = $excel.ActiveWindow
$win.ScrollRow = 22
$win.ScrollColumn = 1
$ws.Range("B49").Select()
